$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Fix the Ethan Virtudazo-style cell (s=3) BEFORE the value reshuffle,
# by copying that distinctive font formatting from its old location (B7) to its
# new location (B17), then restoring B7 to the regular style (copy from B2).
$ws.Range("B7").Copy()
$ws.Range("B17").PasteSpecial(-4122)
$ws.Range("B2").Copy()
$ws.Range("B7").PasteSpecial(-4122)
$ws.Range("B27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 2: Write the new duty-roster names into column B (row 2..31), and
# clear row 32.
$ws.Range("B2").Value = "山口玲, 日高泰聖"
$ws.Range("B3").Value = "志塚惇希"
$ws.Range("B4").Value = "山口洸翔"
$ws.Range("B5").Value = "白岩詩佑介"
$ws.Range("B6").Value = "石井海成, Nicholas Tristan Aryasatyo"
$ws.Range("B7").Value = "小溝賢"
$ws.Range("B8").Value = "小野文哉"
$ws.Range("B9").Value = "渡部魁"
$ws.Range("B10").Value = "崎谷航平, Jun Seomun"
$ws.Range("B11").Value = "三神佳誠"
$ws.Range("B12").Value = "氏家琉貴"
$ws.Range("B13").Value = "羽賀尚生, 島田実"
$ws.Range("B14").Value = "フロアミなので、みな"
$ws.Range("B15").Value = "足立耕平"
$ws.Range("B16").Value = "遠藤隼人"
$ws.Range("B17").Value = "Ethan Virtudazo"
$ws.Range("B18").Value = "富澤天音"
$ws.Range("B19").Value = "神山修造"
$ws.Range("B20").Value = "川田涼介"
$ws.Range("B21").Value = "兒島大志郎"
$ws.Range("B22").Value = "高野怜央"
$ws.Range("B23").Value = "山口玲"
$ws.Range("B24").Value = "志塚惇希"
$ws.Range("B25").Value = "山口洸翔"
$ws.Range("B26").Value = "白岩詩佑介"
$ws.Range("B27").Value = "石井海成"
$ws.Range("B28").Value = "小溝賢"
$ws.Range("B29").Value = "小野文哉"
$ws.Range("B30").Value = "渡部魁"
$ws.Range("B31").Value = "崎谷航平"
$ws.Range("B32").ClearContents()

# Step 3: Update the active selection to match the saved view state.
$ws.Range("E14").Select()
